# Auto-generated edit script for 合肥-漫展信息.xlsx
# Reorders / refreshes the event-listing rows on sheets "展览" and "全部类型"
# to match the scraped data snapshot described in the commit diff.

$wb = $excel.ActiveWorkbook

$dates = @(
  "2024.01.20",
  "2024.01.27",
  "2024.01.28",
  "2024.01.28",
  "2024.01.28",
  "2024.01.29",
  "2024.01.31",
  "2024.02.03",
  "2024.02.04",
  "2024.02.04",
  "2024.02.13",
  "2024.02.14",
  "2024.02.17",
  "2024.02.19",
  "2024.04.04"
)

$names = @(
  "合肥·国乙only新春年会版",
  "合肥·SP同人展·次元派对",
  "合肥·第十二届次元之门动漫游戏博览会-吴磊专场",
  "合肥·第十二届次元之门动漫游戏博览会-赵乾景专场",
  "巢湖·原神&崩铁&崩坏only",
  "肥东· 原神&崩铁&崩坏only",
  "肥西·原神&崩铁&崩坏only",
  "合肥·环形宇宙动漫游戏嘉年华",
  "合肥·环形宇宙动漫游戏嘉年华—吴晛专场",
  "巢湖·原×铁×崩only",
  "合肥·新春AG动漫游戏盛典热血plus",
  "合肥·梦时空SPO1动漫展",
  "合肥·2024运动新春动漫庆典（全ip）",
  "合肥·安徽马娘only",
  "合肥· 第二届漫画城市动漫展 -故事再次开始"
)

$locations = @(
  "文忠路1865号 赫拉诺言艺术中心",
  "临泉路88号板桥里墨园E区1号省羽体中心 省羽体super速搏羽毛球馆",
  "南京路与庐州大道交汇处 合肥滨湖国际会展中心",
  "南京路与庐州大道交汇处 合肥滨湖国际会展中心",
  "团结东路7号 巢湖宾馆",
  "长江东路徽商城2幢B座(祥和地铁站C口步行370米) 曼斯顿尚品酒店",
  "仙满楼·麦肯希酒店 仙满楼·麦肯希酒店",
  "南京路与庐州大道交汇处 合肥滨湖国际会展中心",
  "南京路与庐州大道交汇处 合肥滨湖国际会展中心",
  "健康东路7号 巢湖国际饭店",
  "山西路与太原路交叉口 挥动体育",
  "阜阳路16号 银瑞林国际大酒店",
  "锦绣大道与清潭路交口东北角 李宁体育公园",
  "桐城路与庐江路交叉口西南80米 赤阑桥文玩大厦",
  "凤淮路与固镇路西北角 庐阳全民健身中心"
)

$timeranges = @(
  "2024.01.20 09:30-01.20 17:30",
  "2024.01.27 10:00-01.28 17:00",
  "2024.01.28 10:00-01.28 17:00",
  "2024.01.28 10:00-01.28 17:00",
  "2024.01.28 10:00-01.28 17:00",
  "2024.01.29 10:00-01.29 17:00",
  "2024.01.31 10:00-01.31 17:00",
  "2024.02.03 09:30-02.04 17:00",
  "2024.02.04 11:30-02.04 17:00",
  "2024.02.04 10:00-02.04 17:00",
  "2024.02.13 09:30-02.14 16:00",
  "2024.02.14 10:00-02.14 17:00",
  "2024.02.17 09:00-02.17 17:00",
  "2024.02.19 09:00-02.19 17:00",
  "2024.04.04 09:00-04.05 17:00"
)

$wantcounts = @(
  420,
  1431,
  532,
  315,
  12,
  6,
  11,
  5341,
  135,
  5,
  1703,
  57,
  1062,
  266,
  5474
)

$prices = @(
  "已售罄",
  "55",
  "已售罄",
  "258",
  "55",
  "55",
  "55",
  "65",
  "168",
  "60",
  "39",
  "60",
  "65",
  "68",
  "60"
)

$hasStage = @(
  $false,
  $false,
  $false,
  $false,
  $false,
  $false,
  $false,
  $true,
  $false,
  $false,
  $false,
  $false,
  $false,
  $false,
  $false
)

$links = @(
  "https://show.bilibili.com/platform/detail.html?id=80352&msource=Msearch_colligation",
  "https://show.bilibili.com/platform/detail.html?id=78076&msource=Msearch_colligation",
  "https://show.bilibili.com/platform/detail.html?id=80254&msource=Msearch_colligation",
  "https://show.bilibili.com/platform/detail.html?id=80246&msource=Msearch_colligation",
  "https://show.bilibili.com/platform/detail.html?id=80917&msource=Msearch_colligation",
  "https://show.bilibili.com/platform/detail.html?id=80919&msource=Msearch_colligation",
  "https://show.bilibili.com/platform/detail.html?id=80944&msource=Msearch_colligation",
  "https://show.bilibili.com/platform/detail.html?id=79963&msource=Msearch_colligation",
  "https://show.bilibili.com/platform/detail.html?id=80551&msource=Msearch_colligation",
  "https://show.bilibili.com/platform/detail.html?id=80974&msource=Msearch_colligation",
  "https://show.bilibili.com/platform/detail.html?id=80584&msource=Msearch_colligation",
  "https://show.bilibili.com/platform/detail.html?id=80207&msource=Msearch_colligation",
  "https://show.bilibili.com/platform/detail.html?id=79918&msource=Msearch_colligation",
  "https://show.bilibili.com/platform/detail.html?id=78286&msource=Msearch_colligation",
  "https://show.bilibili.com/platform/detail.html?id=78898&msource=Msearch_colligation"
)

$covers = @(
  "//i2.hdslb.com/bfs/openplatform/202312/VBekVPuH1703840712015.jpeg",
  "//i1.hdslb.com/bfs/openplatform/202311/2v00jbxM1698999146733.jpeg",
  "//i0.hdslb.com/bfs/openplatform/202312/9ClQwbVE1703668101900.jpeg",
  "//i0.hdslb.com/bfs/openplatform/202312/aHzqArm61703662347629.jpeg",
  "//i0.hdslb.com/bfs/openplatform/202401/UekMeUjQ1705462868391.jpeg",
  "//i0.hdslb.com/bfs/openplatform/202401/9XumHIT31705464002179.jpeg",
  "//i0.hdslb.com/bfs/openplatform/202401/euD63Mlp1705479140627.jpeg",
  "//i0.hdslb.com/bfs/openplatform/202312/tBk3WVyX1702968658234.jpeg",
  "//i0.hdslb.com/bfs/openplatform/202401/MSS7qIQp1704695420767.jpeg",
  "//i0.hdslb.com/bfs/openplatform/202401/wVVrdShB1705487994232.jpeg",
  "//i1.hdslb.com/bfs/openplatform/202401/yI94srFk1704703809648.jpeg",
  "//i2.hdslb.com/bfs/openplatform/202312/tQQOHYE01703574162111.jpeg",
  "//i0.hdslb.com/bfs/openplatform/202312/vzuMc0sJ1702902061660.jpeg",
  "//i1.hdslb.com/bfs/openplatform/202311/721L5pIZ1699428443216.jpeg",
  "//i2.hdslb.com/bfs/openplatform/202311/244eBWip1700711342120.jpeg"
)

$rowCount = $dates.Length

foreach ($sheetIndex in 1,4) {
  $ws = $wb.Worksheets.Item($sheetIndex)

  # Drop the last data row (row 17): the refreshed dataset only has 15 events (rows 2-16)
  $ws.Rows.Item(17).Delete()

  # Force text storage for columns that would otherwise be auto-coerced by Excel
  # (dates like "2024.01.20" -> date serials, price codes like "55" -> numbers)
  $ws.Range("B2:B16").NumberFormat = "@"
  $ws.Range("G2:G16").NumberFormat = "@"

  for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 2).Value = $dates[$i]
    $ws.Cells.Item($r, 3).Value = $names[$i]
    $ws.Cells.Item($r, 4).Value = $locations[$i]
    $ws.Cells.Item($r, 5).Value = $timeranges[$i]
    $ws.Cells.Item($r, 6).Value = $wantcounts[$i]
    $ws.Cells.Item($r, 7).Value = $prices[$i]
    $ws.Cells.Item($r, 8).Value = $hasStage[$i]
    $ws.Cells.Item($r, 9).Value = $links[$i]
    $ws.Cells.Item($r, 10).Value = $covers[$i]
  }

  # Restore the original (default) number format now that the text values are locked in
  $ws.Range("B2:B16").NumberFormat = "General"
  $ws.Range("G2:G16").NumberFormat = "General"
}
